# Adds "Web Form" / "Web Form Metadata" / "Web Form Step" record-type
# configuration rows (and their Option Set picklist values) to the
# Customisations workbook, mirroring the existing adx_* entity rows.

$wb = $excel.ActiveWorkbook

$wsRecordTypes = $wb.Worksheets.Item("Record Types")
$wsFields      = $wb.Worksheets.Item("Option Sets")   # placeholder, reset below
$wsFields      = $wb.Worksheets.Item("Fields")
$wsOptionSets  = $wb.Worksheets.Item("Option Sets")

# ---------------------------------------------------------------------
# 1. "Record Types" sheet - three new rows describing the entities
#    adx_webform / adx_webformmetadata / adx_webformstep.
# ---------------------------------------------------------------------
$recordTypeRows = @(
    @("Web Form",          "adx_webform",         "Web Forms",
      "Defines the necessary properties and relationships to the other key entities in order to control the initialization of the form within a web portal."),
    @("Web Form Metadata",  "adx_webformmetadata", "Web Form Metadata",
      "Defines the additional behavior modification logic to augment or override the functionality of form fields that is not possible with entity and form metadata."),
    @("Web Form Step",      "adx_webformstep",     "Web Form Steps",
      "Defines the flow logic of the form's user experience such as steps and conditional branching.")
)

$r = 12
foreach ($row in $recordTypeRows) {
    $wsRecordTypes.Cells.Item($r, 1).Value = $row[0]
    $wsRecordTypes.Cells.Item($r, 2).Value = $row[1]
    $wsRecordTypes.Cells.Item($r, 3).Value = $row[2]
    $wsRecordTypes.Cells.Item($r, 4).Value = $row[3]
    $wsRecordTypes.Cells.Item($r, 5).Value  = $false
    $wsRecordTypes.Cells.Item($r, 6).Value  = $false
    $wsRecordTypes.Cells.Item($r, 7).Value  = $false
    $wsRecordTypes.Cells.Item($r, 8).Value  = $false
    $wsRecordTypes.Cells.Item($r, 9).Value  = $false
    $wsRecordTypes.Cells.Item($r, 10).Value = $false
    $wsRecordTypes.Cells.Item($r, 11).Value = $false
    $r++
}

# ---------------------------------------------------------------------
# 2. "Fields" sheet - field / relationship metadata rows for the new
#    entities (rows 36-44).
# ---------------------------------------------------------------------
function Set-FieldRow {
    param($sheet, $r, $values)

    if ($values.ContainsKey("A")) { $sheet.Cells.Item($r, 1).Value  = $values["A"] }
    if ($values.ContainsKey("B")) { $sheet.Cells.Item($r, 2).Value  = $values["B"] }
    if ($values.ContainsKey("C")) { $sheet.Cells.Item($r, 3).Value  = $values["C"] }
    if ($values.ContainsKey("D")) { $sheet.Cells.Item($r, 4).Value  = $values["D"] }
    if ($values.ContainsKey("E")) { $sheet.Cells.Item($r, 5).Value  = $values["E"] }
    if ($values.ContainsKey("F")) { $sheet.Cells.Item($r, 6).Value  = $values["F"] }
    if ($values.ContainsKey("G")) { $sheet.Cells.Item($r, 7).Value  = $values["G"] }
    if ($values.ContainsKey("H")) { $sheet.Cells.Item($r, 8).Value  = $values["H"] }
    if ($values.ContainsKey("I")) { $sheet.Cells.Item($r, 9).Value  = $values["I"] }
    if ($values.ContainsKey("L")) { $sheet.Cells.Item($r, 12).Value = $values["L"] }
    if ($values.ContainsKey("M")) { $sheet.Cells.Item($r, 13).Value = $values["M"] }
    if ($values.ContainsKey("N")) { $sheet.Cells.Item($r, 14).Value = $values["N"] }
    if ($values.ContainsKey("O")) { $sheet.Cells.Item($r, 15).Value = $values["O"] }
    if ($values.ContainsKey("Q")) { $sheet.Cells.Item($r, 17).Value = $values["Q"] }
    if ($values.ContainsKey("R")) { $sheet.Cells.Item($r, 18).Value = $values["R"] }
    if ($values.ContainsKey("S")) { $sheet.Cells.Item($r, 19).Value = $values["S"] }
    if ($values.ContainsKey("T")) {
        $sheet.Cells.Item($r, 20).Value = $values["T"]
        if ($values.ContainsKey("TBold") -and $values["TBold"]) {
            $sheet.Cells.Item($r, 20).Font.Bold = $true
        }
    }
}

Set-FieldRow $wsFields 36 @{
    A = "adx_webform"; B = "adx_name"; C = "Name"; D = "String";
    E = "Type the name of the custom entity.";
    F = $true; G = $true; H = $true; I = $true;
    L = "N/A"; M = $false; N = 100; O = "Text";
    Q = "N/A"; R = $false; S = -1; T = "N/A"; TBold = $true
}

Set-FieldRow $wsFields 37 @{
    A = "adx_webformstep"; B = "adx_name"; C = "Name"; D = "String";
    E = "Type the name of the custom entity.";
    F = $true; G = $true; H = $true; I = $true;
    L = "N/A"; M = $false; N = 100; O = "Text";
    Q = "N/A"; R = $false; S = -1; T = "N/A"; TBold = $true
}

Set-FieldRow $wsFields 38 @{
    A = "adx_webformstep"; B = "adx_webform"; C = "adx_webform"; D = "Lookup";
    E = "Unique identifier for Web Form associated with Web Form Step.";
    F = $false; G = $true; H = $false; I = $true;
    L = "adx_webform"; M = $true; N = -1;
    Q = "N/A"; R = $false; S = -1; T = "N/A"; TBold = $true
}

Set-FieldRow $wsFields 39 @{
    A = "adx_webformmetadata"; B = "adx_attributelogicalname"; C = "Attribute Logical Name"; D = "String";
    E = "The name of the attribute field to be modified";
    F = $true; G = $false; H = $true; I = $true;
    L = "N/A"; M = $false; N = 100; O = "Text";
    Q = "N/A"; R = $false; S = -1; T = "N/A"; TBold = $true
}

Set-FieldRow $wsFields 40 @{
    A = "adx_webformmetadata"; B = "adx_webformstep"; C = "adx_webformstep"; D = "Lookup";
    E = "Unique identifier for Web Form Step associated with Web Form Metadata.";
    F = $false; G = $true; H = $false; I = $true;
    L = "adx_webformstep"; M = $true; N = -1;
    Q = "N/A"; R = $false; S = -1; T = "N/A"; TBold = $true
}

Set-FieldRow $wsFields 41 @{
    A = "adx_webformmetadata"; B = "adx_type"; C = "Type"; D = "Picklist";
    F = $false; G = $true; H = $true; I = $true;
    L = "N/A"; M = $false; N = -1;
    Q = "N/A"; R = $false; S = -1; T = "WebFormMtType"
}

Set-FieldRow $wsFields 42 @{
    A = "adx_webformmetadata"; B = "adx_tabname"; C = "Tab Name"; D = "String";
    F = $false; G = $false; H = $true; I = $true;
    L = "N/A"; M = $false; N = 200; O = "Text";
    Q = "N/A"; R = $false; S = -1; T = -1
}

Set-FieldRow $wsFields 43 @{
    A = "adx_webformmetadata"; B = "adx_sectionname"; C = "Section Name"; D = "String";
    F = $false; G = $false; H = $true; I = $true;
    L = "N/A"; M = $false; N = 200; O = "Text";
    Q = "N/A"; R = $false; S = -1; T = -1
}

Set-FieldRow $wsFields 44 @{
    A = "adx_webformmetadata"; B = "adx_subgrid_name"; C = "Subgrid Name"; D = "String";
    F = $false; G = $false; H = $true; I = $true;
    L = "N/A"; M = $false; N = 150; O = "Text";
    Q = "N/A"; R = $false; S = -1; T = -1
}

# ---------------------------------------------------------------------
# 3. "Option Sets" sheet - picklist values for the new WebFormMtType
#    option set (rows 11-17).
# ---------------------------------------------------------------------
$optionSetRows = @(
    @(756150000, "Timeline"),
    @(100000000, "Attribute"),
    @(100000005, "Notes"),
    @(100000003, "Purchase"),
    @(100000001, "Section"),
    @(100000004, "Subgrid"),
    @(100000002, "Tab")
)

$r = 11
foreach ($row in $optionSetRows) {
    $wsOptionSets.Cells.Item($r, 1).Value = "WebFormMtType"
    $wsOptionSets.Cells.Item($r, 2).Value = "WebFormMtType"
    $wsOptionSets.Cells.Item($r, 3).Value = $false
    $wsOptionSets.Cells.Item($r, 4).Value = $row[0]
    $wsOptionSets.Cells.Item($r, 5).Value = $row[1]
    $r++
}

# ---------------------------------------------------------------------
# 4. View state: selections on every sheet + which sheet/tab is active.
# ---------------------------------------------------------------------
[void]$wsRecordTypes.Range("E32").Select()
[void]$wsFields.Range("T33").Select()
[void]$wsOptionSets.Range("I29").Select()

$wsOptionSets.Activate()
